# Trade #27 closed at 2026-02-17 20:07:51 - unknown UNKNOWN +0.000%
#
# 1) Summary sheet: bump Total Trades (B6) and recompute Win Rate % (B9)
# 2) Strategy Status sheet: bump MarketMaking row Trades (D5) and Win Rate % (G5)
# 3) All Trades + MarketMaking sheets: append new trade row (#27 -> row 28)

$wb = $excel.ActiveWorkbook

# --- Summary sheet -------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 27
$summary.Range("B9").Value = 44.44

# --- Strategy Status sheet -------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D5").Value = 27
$status.Range("G5").Value = 44.44

# --- New trade row data ----------------------------------------------------
$tradeNum = 27
$date = "2026-02-17"
$time = "20:07:45"
$strategy = "MarketMaking"
$side = "DOWN"
$entryPrice = 0.01
$exitPrice = 0.01
$status_ = "CLOSED"
$pnlPct = 0
$pnlDollar = 0
$capitalAfter = 99.8
$entrySlippage = 0
$exitSlippage = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason = "early_exit"
$duration = 0.11

function Set-TradeRow($ws, $rowIndex) {
    $ws.Cells.Item($rowIndex, 1).Value = $tradeNum
    # Force column B to text so the "YYYY-MM-DD" string isn't auto-converted
    # into a date serial number by Excel's type inference.
    $ws.Cells.Item($rowIndex, 2).NumberFormat = "@"
    $ws.Cells.Item($rowIndex, 2).Value = $date
    $ws.Cells.Item($rowIndex, 3).Value = $time
    $ws.Cells.Item($rowIndex, 4).Value = $strategy
    $ws.Cells.Item($rowIndex, 5).Value = $side
    $ws.Cells.Item($rowIndex, 6).Value = $entryPrice
    $ws.Cells.Item($rowIndex, 7).Value = $exitPrice
    $ws.Cells.Item($rowIndex, 8).Value = $status_
    $ws.Cells.Item($rowIndex, 9).Value = $pnlPct
    $ws.Cells.Item($rowIndex, 10).Value = $pnlDollar
    $ws.Cells.Item($rowIndex, 11).Value = $capitalAfter
    $ws.Cells.Item($rowIndex, 12).Value = $entrySlippage
    $ws.Cells.Item($rowIndex, 13).Value = $exitSlippage
    $ws.Cells.Item($rowIndex, 14).Value = $confidence
    $ws.Cells.Item($rowIndex, 15).Value = $entryReason
    $ws.Cells.Item($rowIndex, 16).Value = $exitReason
    $ws.Cells.Item($rowIndex, 17).Value = $duration
}

# --- All Trades sheet -------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Set-TradeRow $allTrades 28

# --- MarketMaking sheet -------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Set-TradeRow $marketMaking 28
